# Commit: "add recognice of images and click"
# Adds two new FINDWORDANDCREATEAREGION/FINDWORDANDCLICK steps (NSO region) after the
# existing CREATEREGIONAPP,REMOTEDESKTOPNSO step, and appends two new steps at the
# bottom of the script (FINDIAMGEANDCREATEAREGION / CLICKINREGION) for the
# "Aceptar_sentra" region, plus a couple of trailing blank rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. C11 loses its special (now-recycled) style, takes on the plain "wrapped values"
#        style that C1/C4/... already use.
$ws.Range("C1").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# --- 2. Insert a new row at 17 (pushes the old rows 17.. down by one), then populate it.
$ws.Rows.Item(17).Insert(-4121)

# Give the new row the same look & feel as its neighbours (style 0 / style 0-ish),
# then fill in the values.
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C17").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("B17").Value = "FINDWORDANDCREATEAREGION"
$ws.Range("C17").Value = "REMOTEDESKTOPNSO,NSO"

# --- 3. The row that used to be 17 (now 18, "FINDWORDANDCLICK") keeps its keyword but
#        its value becomes NSO,NSO instead of REMOTEDESKTOPNSO,NSO.
$ws.Range("C18").Value = "NSO,NSO"

# --- 4. Rows 19-31 (old 18-30) keep their content unchanged after the shift, nothing to do.

# --- 5. The last row (old 31, now 32) used to read
#        FINDWORDANDCLICK / REMOTEDESKTOPNSO,Aceptar -- replace it with the new
#        "find image, create region" step, and append a new CLICKINREGION step plus two
#        blank trailer rows.
$ws.Range("B32").Value = "FINDIAMGEANDCREATEAREGION"
$ws.Range("C32").Value = "REMOTEDESKTOPNSO,Aceptar_sentra"

# B32 gets a distinctive font (green JetBrains Mono) to flag the new keyword.
$ws.Range("B32").Font.Color = 5867370
$ws.Range("B32").Font.Name = "JetBrains Mono"

# New row 33: CLICKINREGION / Aceptar_sentra,Aceptar_sentra
$ws.Range("B11").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$ws.Range("C33").Value = "Aceptar_sentra,Aceptar_sentra"
$ws.Range("B33").Value = "CLICKINREGION"

# Trailing blank rows 34 & 35 (C only, same style as the other blank-row markers).
$ws.Range("C28").Copy()
$ws.Range("C34").PasteSpecial(-4122)
$ws.Range("C28").Copy()
$ws.Range("C35").PasteSpecial(-4122)
